# Commit: "Removed url in config related to LL feature"
# Adds a new SignIn (LL/Login) test-data row to the SignIn worksheet:
#   SDET185 / @SDET!*% / Vaid username and password

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SignIn")

# New row of test data, right below the existing last row (row 4)
$ws.Range("A5").Value = "SDET185"
$ws.Range("B5").Value = "@SDET!*%"
$ws.Range("C5").Value = "Vaid username and password"

# A5:B5 should look like the other "code-style" (Consolas) cells used
# elsewhere in the workbook, e.g. Register!D9 -- copy that formatting so
# the existing shared style gets reused instead of minting a new one.
$fmtSource = $wb.Worksheets.Item("Register").Range("D9")
$fmtSource.Copy()
$ws.Range("A5:B5").PasteSpecial(-4122)

# C5 keeps a plain Calibri font (distinct cell style from the rest of the
# column, matching the source workbook).
$ws.Range("C5").Font.Name = "Calibri"

# Leave the sheet's selection where the author left it when saving.
$ws.Activate() | Out-Null
$ws.Range("C11").Select() | Out-Null
